# Scene.xlsx: update RelivePos values ("186,6.89,88" -> "186,0,88")
# and move the active cell selection on Sheet1 from E4 to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("RelivePos") rows 2-4 all contained "186,6.89,88";
# change them to "186,0,88".
$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

# Update the saved selection/active cell from E4 to F7.
$ws.Range("F7").Select() | Out-Null
